# Apply the 31-May-2020 COVID data refresh (paises.xlsx) -
# updates the "Pais" sheet with newer totals and re-sorts a few
# countries whose case counts now rank in a different order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country labels for the rows whose relative ranking swapped ---
$ws.Range("A31").Value = "Sudafrica"
$ws.Range("A32").Value = "Portugal"
$ws.Range("A137").Value = "Mauritania"
$ws.Range("A138").Value = "Tanzania"
$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("A201").Value = "Belice"
$ws.Range("A204").Value = "Surinam"
$ws.Range("A205").Value = "San Cristobal y Nieves"
$ws.Range("A210").Value = "Seychelles"
$ws.Range("A211").Value = "Montserrat"
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("A214").Value = "Islas Virgenes Britanicas"

# --- Refreshed figures: Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes ---
$ws.Range("B4").Value = 1828308
$ws.Range("C4").Value = 11488
$ws.Range("D4").Value = 538587
$ws.Range("E4").Value = 1183803
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 361
$ws.Range("H4").Value = 105918
$ws.Range("B10").Value = 190603
$ws.Range("C10").Value = 8776
$ws.Range("D10").Value = 91830
$ws.Range("E10").Value = 93367
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 221
$ws.Range("H10").Value = 5406
$ws.Range("B11").Value = 188882
$ws.Range("C11").Value = 257
$ws.Range("D11").Value = 68355
$ws.Range("E11").Value = 91725
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 31
$ws.Range("H11").Value = 28802
$ws.Range("B12").Value = 183452
$ws.Range("C12").Value = 158
$ws.Range("D12").Value = 165200
$ws.Range("E12").Value = 9650
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = 8602
$ws.Range("B31").Value = 32683
$ws.Range("C31").Value = 1716
$ws.Range("D31").Value = 16809
$ws.Range("E31").Value = 15191
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 40
$ws.Range("H31").Value = 683
$ws.Range("B32").Value = 32500
$ws.Range("C32").Value = 297
$ws.Range("D32").Value = 19409
$ws.Range("E32").Value = 11681
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 14
$ws.Range("H32").Value = 1410
$ws.Range("B79").Value = 3623
$ws.Range("C79").Value = 77
$ws.Range("D79").Value = 2837
$ws.Range("E79").Value = 771
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 15
$ws.Range("B112").Value = 1137
$ws.Range("C112").Value = 15
$ws.Range("D112").Value = 872
$ws.Range("E112").Value = 232
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 33
$ws.Range("B137").Value = 530
$ws.Range("C137").Value = 47
$ws.Range("D137").Value = 27
$ws.Range("E137").Value = 480
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 3
$ws.Range("H137").Value = 23
$ws.Range("B138").Value = 509
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 183
$ws.Range("E138").Value = 305
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 21
$ws.Range("B200").Value = 18
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 18
$ws.Range("E200").Value = 0
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 0
$ws.Range("B201").Value = 18
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 16
$ws.Range("E201").Value = 0
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 2
$ws.Range("B204").Value = 15
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 9
$ws.Range("E204").Value = 0
$ws.Range("F204").Value = 0
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 1
$ws.Range("B205").Value = 15
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 15
$ws.Range("E205").Value = 0
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0
$ws.Range("B210").Value = 11
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 11
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0
$ws.Range("B211").Value = 11
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 10
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 1
$ws.Range("B213").Value = 8
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 8
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0
$ws.Range("B214").Value = 8
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 7
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1

# --- "Last refreshed" banner ---
$ws.Range("A1").Value = "Datos actualizados a 31 de Mayo de 2020 a las 21:05"
